{"js": "// Replace the date line and each three-digit x one-digit multiplication\n// problem's text with its updated value, per the target diff.\nconst replacements = [\n  [\"2025-12-16 Tuesday\", \"2025-12-17 Wednesday\"],\n  [\"521\u00d73=1563\", \"872\u00d78=6976\"],\n  [\"124\u00d74=496\", \"649\u00d74=2596\"],\n  [\"638\u00d75=3190\", \"809\u00d77=5663\"],\n  [\"651\u00d74=2604\", \"182\u00d74=728\"],\n  [\"656\u00d72=1312\", \"984\u00d72=1968\"],\n  [\"970\u00d75=4850\", \"246\u00d74=984\"],\n  [\"827\u00d79=7443\", \"329\u00d77=2303\"],\n  [\"948\u00d79=8532\", \"626\u00d79=5634\"],\n  [\"540\u00d77=3780\", \"868\u00d74=3472\"],\n  [\"579\u00d73=1737\", \"947\u00d77=6629\"],\n  [\"115\u00d77=805\", \"392\u00d78=3136\"],\n  [\"683\u00d72=1366\", \"324\u00d78=2592\"],\n  [\"595\u00d79=5355\", \"577\u00d75=2885\"],\n  [\"618\u00d77=4326\", \"915\u00d74=3660\"],\n  [\"560\u00d72=1120\", \"743\u00d77=5201\"],\n  [\"375\u00d79=3375\", \"873\u00d78=6984\"],\n  [\"691\u00d75=3455\", \"605\u00d79=5445\"],\n  [\"858\u00d75=4290\", \"716\u00d72=1432\"],\n  [\"312\u00d74=1248\", \"321\u00d79=2889\"],\n  [\"807\u00d77=5649\", \"815\u00d79=7335\"],\n  [\"889\u00d79=8001\", \"164\u00d74=656\"],\n  [\"666\u00d78=5328\", \"381\u00d78=3048\"],\n  [\"473\u00d78=3784\", \"690\u00d75=3450\"],\n  [\"716\u00d78=5728\", \"973\u00d77=6811\"],\n  [\"649\u00d73=1947\", \"657\u00d78=5256\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each three-digit x one-digit multiplication\n# problem's text with its updated value, per the target diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-12-16 Tuesday\", \"2025-12-17 Wednesday\"),\n    @(\"521\u00d73=1563\", \"872\u00d78=6976\"),\n    @(\"124\u00d74=496\", \"649\u00d74=2596\"),\n    @(\"638\u00d75=3190\", \"809\u00d77=5663\"),\n    @(\"651\u00d74=2604\", \"182\u00d74=728\"),\n    @(\"656\u00d72=1312\", \"984\u00d72=1968\"),\n    @(\"970\u00d75=4850\", \"246\u00d74=984\"),\n    @(\"827\u00d79=7443\", \"329\u00d77=2303\"),\n    @(\"948\u00d79=8532\", \"626\u00d79=5634\"),\n    @(\"540\u00d77=3780\", \"868\u00d74=3472\"),\n    @(\"579\u00d73=1737\", \"947\u00d77=6629\"),\n    @(\"115\u00d77=805\", \"392\u00d78=3136\"),\n    @(\"683\u00d72=1366\", \"324\u00d78=2592\"),\n    @(\"595\u00d79=5355\", \"577\u00d75=2885\"),\n    @(\"618\u00d77=4326\", \"915\u00d74=3660\"),\n    @(\"560\u00d72=1120\", \"743\u00d77=5201\"),\n    @(\"375\u00d79=3375\", \"873\u00d78=6984\"),\n    @(\"691\u00d75=3455\", \"605\u00d79=5445\"),\n    @(\"858\u00d75=4290\", \"716\u00d72=1432\"),\n    @(\"312\u00d74=1248\", \"321\u00d79=2889\"),\n    @(\"807\u00d77=5649\", \"815\u00d79=7335\"),\n    @(\"889\u00d79=8001\", \"164\u00d74=656\"),\n    @(\"666\u00d78=5328\", \"381\u00d78=3048\"),\n    @(\"473\u00d78=3784\", \"690\u00d75=3450\"),\n    @(\"716\u00d78=5728\", \"973\u00d77=6811\"),\n    @(\"649\u00d73=1947\", \"657\u00d78=5256\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Replacement.ClearFormatting()\n    $r.Find.Text = $oldText\n    $r.Find.Replacement.Text = $newText\n    $r.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\n$d.Save()\n"}
